# Auto-generated Excel COM-interop script applying scheduled-runner data refresh
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 449
$ws.Cells.Item(9, 9).Value = 469.75
$ws.Cells.Item(9, 11).Value = 469.75
$ws.Cells.Item(9, 13).Value = -300.75
$ws.Cells.Item(33, 8).Value = 29
$ws.Cells.Item(33, 9).Value = 29
$ws.Cells.Item(33, 11).Value = 29
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(74, 8).Value = 138476.64
$ws.Cells.Item(74, 9).Value = 148029.8
$ws.Cells.Item(74, 10).Value = 42945
$ws.Cells.Item(74, 11).Value = 148029.8
$ws.Cells.Item(74, 12).Value = 42945
$ws.Cells.Item(74, 13).Value = -147093.8
$ws.Cells.Item(74, 14).Value = -44817
$ws.Cells.Item(77, 8).Value = 138476.64
$ws.Cells.Item(77, 9).Value = 148029.8
$ws.Cells.Item(77, 10).Value = 42945
$ws.Cells.Item(77, 11).Value = 740149
$ws.Cells.Item(77, 12).Value = 214725
$ws.Cells.Item(77, 13).Value = -735469
$ws.Cells.Item(77, 14).Value = -224085
$ws.Cells.Item(98, 8).Value = 1961.6666
$ws.Cells.Item(98, 10).Value = 2222
$ws.Cells.Item(98, 12).Value = 2222
$ws.Cells.Item(98, 14).Value = -5218
$ws.Cells.Item(115, 8).Value = 799
$ws.Cells.Item(115, 9).Value = 799
$ws.Cells.Item(115, 11).Value = 2397
$ws.Cells.Item(115, 13).Value = -830
$ws.Cells.Item(122, 8).Value = 1961.6666
$ws.Cells.Item(122, 10).Value = 2222
$ws.Cells.Item(122, 12).Value = 6666
$ws.Cells.Item(122, 14).Value = -11566
$ws.Cells.Item(138, 8).Value = 1960.6875
$ws.Cells.Item(138, 9).Value = 1960.6875
$ws.Cells.Item(138, 11).Value = 5882.0625
$ws.Cells.Item(138, 13).Value = -742.0625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 10).Value = 2333997
$ws.Cells.Item(32, 12).Value = 2333997
$ws.Cells.Item(32, 14).Value = -2334571
$ws.Cells.Item(74, 8).Value = 1729.4286
$ws.Cells.Item(74, 9).Value = 1702.75
$ws.Cells.Item(74, 10).Value = 1765
$ws.Cells.Item(74, 11).Value = 1702.75
$ws.Cells.Item(74, 12).Value = 1765
$ws.Cells.Item(74, 13).Value = -828.75
$ws.Cells.Item(74, 14).Value = -3513
$ws.Cells.Item(77, 8).Value = 1729.4286
$ws.Cells.Item(77, 9).Value = 1702.75
$ws.Cells.Item(77, 10).Value = 1765
$ws.Cells.Item(77, 11).Value = 8513.75
$ws.Cells.Item(77, 12).Value = 8825
$ws.Cells.Item(77, 13).Value = -4145.75
$ws.Cells.Item(77, 14).Value = -17561
$ws.Cells.Item(92, 8).Value = 56750
$ws.Cells.Item(92, 10).Value = 56750
$ws.Cells.Item(92, 12).Value = 56750
$ws.Cells.Item(92, 14).Value = -61742
$ws.Cells.Item(95, 8).Value = 8000
$ws.Cells.Item(95, 10).Value = 8000
$ws.Cells.Item(95, 12).Value = 8000
$ws.Cells.Item(95, 14).Value = -13492
$ws.Cells.Item(106, 8).Value = 18500
$ws.Cells.Item(106, 10).Value = 18500
$ws.Cells.Item(106, 12).Value = 18500
$ws.Cells.Item(106, 14).Value = -21024
$ws.Cells.Item(109, 8).Value = 99999
$ws.Cells.Item(109, 10).Value = 99999
$ws.Cells.Item(109, 12).Value = 99999
$ws.Cells.Item(109, 14).Value = -102773
$ws.Cells.Item(112, 8).Value = 150000
$ws.Cells.Item(112, 10).Value = 150000
$ws.Cells.Item(112, 12).Value = 150000
$ws.Cells.Item(112, 14).Value = -152954
$ws.Cells.Item(114, 8).Value = 150000
$ws.Cells.Item(114, 10).Value = 150000
$ws.Cells.Item(114, 12).Value = 150000
$ws.Cells.Item(114, 14).Value = -158678
$ws.Cells.Item(132, 8).Value = 1388.4286

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 20435.5
$ws.Cells.Item(33, 9).Value = 14621
$ws.Cells.Item(33, 11).Value = 14621
$ws.Cells.Item(33, 13).Value = -14285
$ws.Cells.Item(94, 8).Value = 1578.826
$ws.Cells.Item(94, 9).Value = 1443.5238
$ws.Cells.Item(94, 11).Value = 1443.5238
$ws.Cells.Item(94, 13).Value = -992.5237999999999
$ws.Cells.Item(99, 8).Value = 3000
$ws.Cells.Item(99, 9).Value = 3000
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 3000
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -1502
$ws.Cells.Item(99, 14).ClearContents() | Out-Null
$ws.Cells.Item(134, 8).Value = 1239.5
$ws.Cells.Item(134, 9).Value = 1239.5
$ws.Cells.Item(134, 11).Value = 3718.5
$ws.Cells.Item(134, 13).Value = -1183.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 4133.3335
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 4133.3335
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 4133.3335
$ws.Cells.Item(12, 13).ClearContents() | Out-Null
$ws.Cells.Item(12, 14).Value = -4473.3335
$ws.Cells.Item(31, 8).Value = 2076.8
$ws.Cells.Item(31, 9).Value = 2071
$ws.Cells.Item(31, 10).Value = 2100
$ws.Cells.Item(31, 11).Value = 2071
$ws.Cells.Item(31, 12).Value = 2100
$ws.Cells.Item(31, 13).Value = -1776
$ws.Cells.Item(31, 14).Value = -2690
$ws.Cells.Item(34, 8).Value = 2076.8
$ws.Cells.Item(34, 9).Value = 2071
$ws.Cells.Item(34, 10).Value = 2100
$ws.Cells.Item(34, 11).Value = 2071
$ws.Cells.Item(34, 12).Value = 2100
$ws.Cells.Item(34, 13).Value = -1869
$ws.Cells.Item(34, 14).Value = -2504
$ws.Cells.Item(42, 8).Value = 49500
$ws.Cells.Item(42, 9).Value = 49000
$ws.Cells.Item(42, 10).Value = 50000
$ws.Cells.Item(42, 11).Value = 49000
$ws.Cells.Item(42, 12).Value = 50000
$ws.Cells.Item(42, 13).Value = -48407
$ws.Cells.Item(42, 14).Value = -51186
$ws.Cells.Item(64, 8).Value = 15000
$ws.Cells.Item(64, 10).Value = 15000
$ws.Cells.Item(64, 12).Value = 15000
$ws.Cells.Item(64, 14).Value = -15496
$ws.Cells.Item(67, 8).Value = 15000
$ws.Cells.Item(67, 10).Value = 15000
$ws.Cells.Item(67, 12).Value = 15000
$ws.Cells.Item(67, 14).Value = -16716
$ws.Cells.Item(88, 8).Value = 15999.5
$ws.Cells.Item(88, 9).Value = 12000
$ws.Cells.Item(88, 10).Value = 19999
$ws.Cells.Item(88, 11).Value = 12000
$ws.Cells.Item(88, 12).Value = 19999
$ws.Cells.Item(88, 13).Value = -11594
$ws.Cells.Item(88, 14).Value = -20811
$ws.Cells.Item(91, 8).Value = 15999.5
$ws.Cells.Item(91, 9).Value = 12000
$ws.Cells.Item(91, 10).Value = 19999
$ws.Cells.Item(91, 11).Value = 12000
$ws.Cells.Item(91, 12).Value = 19999
$ws.Cells.Item(91, 13).Value = -10596
$ws.Cells.Item(91, 14).Value = -22807
$ws.Cells.Item(132, 8).Value = 10632.4
$ws.Cells.Item(132, 9).Value = 10632.4
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 31897.2
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -29367.2
$ws.Cells.Item(132, 14).ClearContents() | Out-Null
$ws.Cells.Item(134, 8).Value = 2322.158
$ws.Cells.Item(134, 9).Value = 2242.4119
$ws.Cells.Item(134, 11).Value = 6727.2357
$ws.Cells.Item(134, 13).Value = -4192.2357

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).ClearContents() | Out-Null
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).ClearContents() | Out-Null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2308.4443
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents() | Out-Null
$ws.Cells.Item(122, 8).Value = 4063.1667
$ws.Cells.Item(122, 9).Value = 4063.1667
$ws.Cells.Item(122, 11).Value = 12189.5001
$ws.Cells.Item(122, 13).Value = -9739.500100000001
$ws.Cells.Item(132, 8).Value = 4937.7
$ws.Cells.Item(132, 9).Value = 5195.9443
$ws.Cells.Item(132, 10).Value = 2613.5
$ws.Cells.Item(132, 11).Value = 15587.8329
$ws.Cells.Item(132, 12).Value = 7840.5
$ws.Cells.Item(132, 13).Value = -13057.8329
$ws.Cells.Item(132, 14).Value = -12900.5
$ws.Cells.Item(134, 8).Value = 50000
$ws.Cells.Item(134, 10).Value = 50000
$ws.Cells.Item(134, 12).Value = 150000
$ws.Cells.Item(134, 14).Value = -155070
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents() | Out-Null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(26, 8).Value = 7000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents() | Out-Null
$ws.Cells.Item(40, 8).Value = 3752.1875
$ws.Cells.Item(40, 9).Value = 3476.2727
$ws.Cells.Item(40, 11).Value = 3476.2727
$ws.Cells.Item(40, 13).Value = -3340.2727
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 13).ClearContents() | Out-Null
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 13).ClearContents() | Out-Null
$ws.Cells.Item(100, 8).Value = 4163.067
$ws.Cells.Item(100, 9).Value = 4372.923
$ws.Cells.Item(100, 10).Value = 2799
$ws.Cells.Item(100, 11).Value = 4372.923
$ws.Cells.Item(100, 12).Value = 2799
$ws.Cells.Item(100, 13).Value = -3831.923
$ws.Cells.Item(100, 14).Value = -3881
$ws.Cells.Item(132, 8).Value = 3141.2
$ws.Cells.Item(132, 9).Value = 3141.2
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 9423.599999999999
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -6893.599999999999
$ws.Cells.Item(132, 14).ClearContents() | Out-Null
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents() | Out-Null
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents() | Out-Null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents() | Out-Null
$ws.Cells.Item(122, 8).Value = 7988.625
$ws.Cells.Item(122, 9).Value = 7514.2
$ws.Cells.Item(122, 11).Value = 22542.6
$ws.Cells.Item(122, 13).Value = -20092.6
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents() | Out-Null
$ws.Cells.Item(136, 8).Value = 1737.75
$ws.Cells.Item(136, 9).Value = 1737.75
$ws.Cells.Item(136, 11).Value = 5213.25
$ws.Cells.Item(136, 13).Value = -2663.25
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).ClearContents() | Out-Null
$ws.Cells.Item(139, 14).ClearContents() | Out-Null
